# Insert two new rows for a new weekly data entry, shifting the
# existing data (rows 375-494) down to rows 377-496.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A375:A376").EntireRow.Insert()

# Populate the newly inserted rows with the new week's data. Most of
# the descriptive columns (Mercado, Region, Categoria, etc.) are the
# same as the row that used to occupy this slot; only the date and the
# price/volume figures are new.
$ws.Cells.Item(375, 1).Value = 8
$ws.Cells.Item(375, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(375, 3).Value = "Coquimbo"
$ws.Cells.Item(375, 4).Value = 45146
$ws.Cells.Item(375, 5).Value = 4
$ws.Cells.Item(375, 6).Value = 100114014
$ws.Cells.Item(375, 7).Value = "Betarraga"
$ws.Cells.Item(375, 8).Value = "Sin especificar"
$ws.Cells.Item(375, 9).Value = "Primera"
$ws.Cells.Item(375, 10).Value = 1900
$ws.Cells.Item(375, 11).Value = 550
$ws.Cells.Item(375, 12).Value = 600
$ws.Cells.Item(375, 13).Value = 575
$ws.Cells.Item(375, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(375, 15).Value = "Provincia del Elqu$([char]0x00ED)"
$ws.Cells.Item(375, 16).Value = 192
$ws.Cells.Item(375, 17).Value = 3
$ws.Cells.Item(375, 18).Value = "Hortaliza"

$ws.Cells.Item(376, 1).Value = 8
$ws.Cells.Item(376, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(376, 3).Value = "Coquimbo"
$ws.Cells.Item(376, 4).Value = 45146
$ws.Cells.Item(376, 5).Value = 4
$ws.Cells.Item(376, 6).Value = 100114014
$ws.Cells.Item(376, 7).Value = "Betarraga"
$ws.Cells.Item(376, 8).Value = "Sin especificar"
$ws.Cells.Item(376, 9).Value = "Segunda"
$ws.Cells.Item(376, 10).Value = 960
$ws.Cells.Item(376, 11).Value = 450
$ws.Cells.Item(376, 12).Value = 500
$ws.Cells.Item(376, 13).Value = 475
$ws.Cells.Item(376, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(376, 15).Value = "Provincia del Elqu$([char]0x00ED)"
$ws.Cells.Item(376, 16).Value = 158
$ws.Cells.Item(376, 17).Value = 3
$ws.Cells.Item(376, 18).Value = "Hortaliza"
